# Update column G ("K") values on Sheet1 for rows 2-27.
# These are the new strikeout (K) counts replacing the old "Strike#" based values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newValues = @{
    2  = 1
    3  = 0
    4  = 0
    5  = 1
    6  = 2
    7  = 3
    8  = 1
    9  = 2
    10 = 2
    11 = 3
    12 = 2
    13 = 2
    14 = 2
    15 = 1
    16 = 1
    17 = 2
    18 = 3
    19 = 0
    20 = 2
    21 = 1
    22 = 3
    23 = 3
    24 = 1
    25 = 1
    26 = 1
    27 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $newValues[$row]
}
